$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("Input")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsRepay  = $wb.Worksheets.Item("Repayment schedule")
$wsTrans  = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------
# Summary sheet updates
# ---------------------------------------------------------------
$wsSummary.Range("A4").Value = 100
$wsSummary.Range("B4").Value = 50
$wsSummary.Range("E4").Value = 50

# Introduce a blank, unstyled G2 cell (mirrors an existing unstyled
# blank cell elsewhere in the workbook so no new cell style is minted).
[void]$wsTrans.Range("K3").Copy()
[void]$wsSummary.Range("G2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Repayment schedule sheet updates
# ---------------------------------------------------------------
$wsRepay.Range("I3").Value = 50
$wsRepay.Range("K3").Value = 937.72
$wsRepay.Range("L3").Value = 937.72

$wsRepay.Range("I5").Value = 50
$wsRepay.Range("K5").Value = 937.72
$wsRepay.Range("P5").Value = 937.72

$wsRepay.Range("P2").Clear()
$wsRepay.Range("O3").Clear()
$wsRepay.Range("O4").Clear()
$wsRepay.Range("O5").Clear()
$wsRepay.Range("O6").Clear()
$wsRepay.Range("O7").Clear()
$wsRepay.Range("O8").Clear()

# ---------------------------------------------------------------
# Transactions sheet updates
# ---------------------------------------------------------------
$wsTrans.Range("A2").Value = 6348
$wsTrans.Range("E2").Value = 937.72
$wsTrans.Range("H2").Value = 50
$wsTrans.Range("A3").Value = 691

# ---------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------
[void]$wsRepay.Activate()
[void]$wsRepay.Range("F6").Select()

[void]$wsTrans.Activate()
[void]$wsTrans.Range("D3").Select()
